$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 1).Value = "红 宝 丽"
$ws.Cells.Item(2, 2).Value = "特变电工"
$ws.Cells.Item(2, 3).Value = "特变电工"
$ws.Cells.Item(3, 1).Value = "特变电工"
$ws.Cells.Item(3, 2).Value = "红 宝 丽"
$ws.Cells.Item(3, 3).Value = "红宝丽"
$ws.Cells.Item(4, 1).Value = "利欧股份"
$ws.Cells.Item(4, 2).Value = "东方财富"
$ws.Cells.Item(4, 3).Value = "利欧股份"
$ws.Cells.Item(5, 1).Value = "蓝色光标"
$ws.Cells.Item(5, 2).Value = "亨通光电"
$ws.Cells.Item(5, 3).Value = "白银有色"
$ws.Cells.Item(6, 1).Value = "天地在线"
$ws.Cells.Item(6, 2).Value = "通鼎互联"
$ws.Cells.Item(6, 3).Value = "蓝色光标"
$ws.Cells.Item(7, 1).Value = "西部材料"
$ws.Cells.Item(7, 2).Value = "西部材料"
$ws.Cells.Item(7, 3).Value = "锋龙股份"
$ws.Cells.Item(8, 1).Value = "信维通信"
$ws.Cells.Item(8, 2).Value = "利欧股份"
$ws.Cells.Item(8, 3).Value = "天奇股份"
$ws.Cells.Item(9, 1).Value = "白银有色"
$ws.Cells.Item(9, 2).Value = "贵州茅台"
$ws.Cells.Item(9, 3).Value = "天地在线"
$ws.Cells.Item(10, 1).Value = "通鼎互联"
$ws.Cells.Item(10, 2).Value = "湖南黄金"
$ws.Cells.Item(10, 3).Value = "航天发展"
$ws.Cells.Item(11, 1).Value = "亨通光电"
$ws.Cells.Item(11, 2).Value = "信维通信"
$ws.Cells.Item(11, 3).Value = "杰瑞股份"
$ws.Cells.Item(12, 1).Value = "湖南黄金"
$ws.Cells.Item(12, 2).Value = "中国西电"
$ws.Cells.Item(12, 3).Value = "浙文互联"
$ws.Cells.Item(13, 1).Value = "浙文互联"
$ws.Cells.Item(13, 2).Value = "蓝色光标"
$ws.Cells.Item(13, 3).Value = "湖南白银"
$ws.Cells.Item(14, 1).Value = "天奇股份"
$ws.Cells.Item(14, 2).Value = "白银有色"
$ws.Cells.Item(14, 3).Value = "湖南黄金"
$ws.Cells.Item(15, 1).Value = "锋龙股份"
$ws.Cells.Item(15, 2).Value = "天奇股份"
$ws.Cells.Item(15, 3).Value = "农发种业"
$ws.Cells.Item(16, 1).Value = "中国西电"
$ws.Cells.Item(16, 2).Value = "农发种业"
$ws.Cells.Item(16, 3).Value = "通鼎互联"
$ws.Cells.Item(17, 1).Value = "东方财富"
$ws.Cells.Item(17, 2).Value = "太极实业"
$ws.Cells.Item(17, 3).Value = "嘉美包装"
$ws.Cells.Item(18, 1).Value = "航天发展"
$ws.Cells.Item(18, 2).Value = "锋龙股份"
$ws.Cells.Item(18, 3).Value = "亨通光电"
$ws.Cells.Item(19, 1).Value = "农发种业"
$ws.Cells.Item(19, 2).Value = "百川股份"
$ws.Cells.Item(19, 3).Value = "中际旭创"
$ws.Cells.Item(20, 1).Value = "保变电气"
$ws.Cells.Item(20, 2).Value = "保变电气"
$ws.Cells.Item(20, 3).Value = "紫金矿业"
$ws.Cells.Item(21, 1).Value = "贵州茅台"
$ws.Cells.Item(21, 2).Value = "天地在线"
$ws.Cells.Item(21, 3).Value = "洲际油气"
